# Documents the new notify() Artisan Command on the "Commands" sheet of the
# Event Custom Buttons workbook, by inserting a new table row right above the
# existing "setCanvasColor(...)" row and filling in the command name and its
# description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Make sure the Commands sheet is (and stays) the active one, as it was
# before the edit.
$ws.Activate()

# Push row 89 ("setCanvasColor(<color>)" ...) and everything below it down
# by one row to make room for the new "notify(...)" entry.
$ws.Rows.Item(89).Insert()

$ws.Cells.Item(89, 2).Value = "notify(<title>,[<msg>])"
$ws.Cells.Item(89, 3).Value = "sends notification with title <title> and optional message <msg>"

# Reflect the edit location as the current selection, like in the authored
# workbook.
$ws.Range("B89:C89").Select()
